# Update "想去人数" (want-to-go count) figures on the "展览", "演出" and
# "全部类型" worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览" = @{
        3  = 12247
        9  = 2615
        10 = 1136
        11 = 211
        12 = 81
        13 = 5353
        15 = 215
        16 = 563
        17 = 11503
        18 = 11616
        23 = 57
    }
    "演出" = @{
        2 = 5
    }
    "全部类型" = @{
        3  = 12247
        9  = 2615
        10 = 5
        11 = 1136
        12 = 211
        13 = 81
        14 = 5353
        16 = 215
        17 = 563
        18 = 11503
        19 = 11616
        24 = 57
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
